$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 830.2941  # H15: was 510.3441
$ws.Cells.Item(15, 9).Value = 830.2941  # I15: was 510.3441
$ws.Cells.Item(15, 11).Value = 2490.8823  # K15: was 1531.0323
$ws.Cells.Item(15, 13).Value = -2321.8823  # M15: was -1362.0323

$ws.Cells.Item(17, 8).Value = 1132.5555  # H17: was 1098.125
$ws.Cells.Item(17, 10).Value = 1028.5927  # J17: was 997.5806
$ws.Cells.Item(17, 12).Value = 3085.7781  # L17: was 2992.7418
$ws.Cells.Item(17, 14).Value = -3421.7781  # N17: was -3328.7418

$ws.Cells.Item(39, 8).Value = 303.375  # H39: was 400.6842
$ws.Cells.Item(39, 9).Value = 113.25  # I39: was 146.4
$ws.Cells.Item(39, 10).Value = 493.5  # J39: was 683.2222
$ws.Cells.Item(39, 11).Value = 339.75  # K39: was 439.2
$ws.Cells.Item(39, 12).Value = 1480.5  # L39: was 2049.6666
$ws.Cells.Item(39, 13).Value = -43.75  # M39: was -143.2
$ws.Cells.Item(39, 14).Value = -2072.5  # N39: was -2641.6666

$ws.Cells.Item(62, 8).Value = 1331  # H62: was 1375.7142
$ws.Cells.Item(62, 9).Value = 1273.2222  # I62: was 1332.5
$ws.Cells.Item(62, 10).Value = 1417.6666  # J62: was 1433.3334
$ws.Cells.Item(62, 11).Value = 1273.2222  # K62: was 1332.5
$ws.Cells.Item(62, 12).Value = 1417.6666  # L62: was 1433.3334
$ws.Cells.Item(62, 13).Value = -649.2221999999999  # M62: was -708.5
$ws.Cells.Item(62, 14).Value = -2665.6666  # N62: was -2681.3334

$ws.Cells.Item(65, 8).Value = 1331  # H65: was 1375.7142
$ws.Cells.Item(65, 9).Value = 1273.2222  # I65: was 1332.5
$ws.Cells.Item(65, 10).Value = 1417.6666  # J65: was 1433.3334
$ws.Cells.Item(65, 11).Value = 6366.111  # K65: was 6662.5
$ws.Cells.Item(65, 12).Value = 7088.333000000001  # L65: was 7166.666999999999
$ws.Cells.Item(65, 13).Value = -3246.111  # M65: was -3542.5
$ws.Cells.Item(65, 14).Value = -13328.333  # N65: was -13406.667

$ws.Cells.Item(138, 8).Value = 3638.7307  # H138: was 3396.8035
$ws.Cells.Item(138, 9).Value = 3219.4  # I138: was 2025.5454
$ws.Cells.Item(138, 10).Value = 3683.3403  # J138: was 3732
$ws.Cells.Item(138, 11).Value = 9658.200000000001  # K138: was 6076.6362
$ws.Cells.Item(138, 12).Value = 11050.0209  # L138: was 11196
$ws.Cells.Item(138, 13).Value = -4518.200000000001  # M138: was -936.6361999999999
$ws.Cells.Item(138, 14).Value = -21330.0209  # N138: was -21476

$ws.Cells.Item(141, 8).Value = 17443.715  # H141: was 20351
$ws.Cells.Item(141, 9).Value = 28000.25  # I141: was 37000.332
$ws.Cells.Item(141, 10).Value = 3368.3333  # J141: was 3701.6667
$ws.Cells.Item(141, 11).Value = 84000.75  # K141: was 111000.996
$ws.Cells.Item(141, 12).Value = 10104.9999  # L141: was 11105.0001
$ws.Cells.Item(141, 13).Value = -78820.75  # M141: was -105820.996
$ws.Cells.Item(141, 14).Value = -20464.9999  # N141: was -21465.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 5312.875  # H31: was 13352.167
$ws.Cells.Item(31, 9).Value = 5312.875  # I31: was 13352.167
$ws.Cells.Item(31, 11).Value = 5312.875  # K31: was 13352.167
$ws.Cells.Item(31, 13).Value = -5018.875  # M31: was -13058.167

$ws.Cells.Item(32, 8).Value = 7143.5967  # H32: was 7918.9243
$ws.Cells.Item(32, 9).Value = 5126.7104  # I32: was 6030.3667
$ws.Cells.Item(32, 10).Value = 10337  # J32: was 10382.261
$ws.Cells.Item(32, 11).Value = 5126.7104  # K32: was 6030.3667
$ws.Cells.Item(32, 12).Value = 10337  # L32: was 10382.261
$ws.Cells.Item(32, 13).Value = -4839.7104  # M32: was -5743.3667
$ws.Cells.Item(32, 14).Value = -10911  # N32: was -10956.261

$ws.Cells.Item(51, 8).Value = 47801  # H51: was 55000
$ws.Cells.Item(51, 10).Value = 47801  # J51: was 55000
$ws.Cells.Item(51, 12).Value = 47801  # L51: was 55000
$ws.Cells.Item(51, 14).Value = -49313  # N51: was -56512

$ws.Cells.Item(119, 8).Value = 0  # H119: was 35488.6
$ws.Cells.Item(119, 10).Value = 0  # J119: was 35488.6
$ws.Cells.Item(119, 12).Value = 0  # L119: was 35488.6
$ws.Cells.Item(119, 14).ClearContents()  # N119: was -45164.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 118842  # H59: was 118844
$ws.Cells.Item(59, 10).Value = 118842  # J59: was 118844
$ws.Cells.Item(59, 12).Value = 118842  # L59: was 118844
$ws.Cells.Item(59, 14).Value = -120536  # N59: was -120538

$ws.Cells.Item(134, 8).Value = 2368.1272  # H134: was 2554.1836
$ws.Cells.Item(134, 9).Value = 1358.2439  # I134: was 1450.2778
$ws.Cells.Item(134, 10).Value = 5325.643  # J134: was 5611.154
$ws.Cells.Item(134, 11).Value = 4074.7317  # K134: was 4350.8334
$ws.Cells.Item(134, 12).Value = 15976.929  # L134: was 16833.462
$ws.Cells.Item(134, 13).Value = -1539.7317  # M134: was -1815.8334
$ws.Cells.Item(134, 14).Value = -21046.929  # N134: was -21903.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 25004586  # H31: was 23814152
$ws.Cells.Item(31, 9).Value = 2470.7  # I31: was 2478.7778
$ws.Cells.Item(31, 10).Value = 50006700  # J31: was 41672908
$ws.Cells.Item(31, 11).Value = 2470.7  # K31: was 2478.7778
$ws.Cells.Item(31, 12).Value = 50006700  # L31: was 41672908
$ws.Cells.Item(31, 13).Value = -2175.7  # M31: was -2183.7778
$ws.Cells.Item(31, 14).Value = -50007290  # N31: was -41673498

$ws.Cells.Item(34, 8).Value = 25004586  # H34: was 23814152
$ws.Cells.Item(34, 9).Value = 2470.7  # I34: was 2478.7778
$ws.Cells.Item(34, 10).Value = 50006700  # J34: was 41672908
$ws.Cells.Item(34, 11).Value = 2470.7  # K34: was 2478.7778
$ws.Cells.Item(34, 12).Value = 50006700  # L34: was 41672908
$ws.Cells.Item(34, 13).Value = -2268.7  # M34: was -2276.7778
$ws.Cells.Item(34, 14).Value = -50007104  # N34: was -41673312

$ws.Cells.Item(99, 8).Value = 12503844  # H99: was 14289966
$ws.Cells.Item(99, 9).Value = 33336318  # I99: was 33336320
$ws.Cells.Item(99, 10).Value = 4359  # J99: was 5200
$ws.Cells.Item(99, 11).Value = 33336318  # K99: was 33336320
$ws.Cells.Item(99, 12).Value = 4359  # L99: was 5200
$ws.Cells.Item(99, 13).Value = -33334820  # M99: was -33334822
$ws.Cells.Item(99, 14).Value = -7355  # N99: was -8196

$ws.Cells.Item(126, 8).Value = 12503844  # H126: was 14289966
$ws.Cells.Item(126, 9).Value = 33336318  # I126: was 33336320
$ws.Cells.Item(126, 10).Value = 4359  # J126: was 5200
$ws.Cells.Item(126, 11).Value = 100008954  # K126: was 100008960
$ws.Cells.Item(126, 12).Value = 13077  # L126: was 15600
$ws.Cells.Item(126, 13).Value = -100006484  # M126: was -100006490
$ws.Cells.Item(126, 14).Value = -18017  # N126: was -20540

$ws.Cells.Item(132, 8).Value = 3755  # H132: was 3656.6667
$ws.Cells.Item(132, 9).Value = 3363.7778  # I132: was 3508.7058
$ws.Cells.Item(132, 10).Value = 4635.25  # J132: was 3908.2
$ws.Cells.Item(132, 11).Value = 10091.3334  # K132: was 10526.1174
$ws.Cells.Item(132, 12).Value = 13905.75  # L132: was 11724.6
$ws.Cells.Item(132, 13).Value = -7561.3334  # M132: was -7996.117400000001
$ws.Cells.Item(132, 14).Value = -18965.75  # N132: was -16784.6

$ws.Cells.Item(134, 8).Value = 7105.909  # H134: was 8486.111000000001
$ws.Cells.Item(134, 9).Value = 15188  # I134: was 21007.2
$ws.Cells.Item(134, 10).Value = 3334.2666  # J134: was 3670.3076
$ws.Cells.Item(134, 11).Value = 45564  # K134: was 63021.60000000001
$ws.Cells.Item(134, 12).Value = 10002.7998  # L134: was 11010.9228
$ws.Cells.Item(134, 13).Value = -43029  # M134: was -60486.60000000001
$ws.Cells.Item(134, 14).Value = -15072.7998  # N134: was -16080.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 1604.762  # H26: was 1838.8889
$ws.Cells.Item(26, 10).Value = 2056.25  # J26: was 2484.6155
$ws.Cells.Item(26, 12).Value = 6168.75  # L26: was 7453.8465
$ws.Cells.Item(26, 14).Value = -6744.75  # N26: was -8029.8465

$ws.Cells.Item(56, 8).Value = 8892.5  # H56: was 10036.667
$ws.Cells.Item(56, 9).Value = 8892.5  # I56: was 10036.667
$ws.Cells.Item(56, 11).Value = 8892.5  # K56: was 10036.667
$ws.Cells.Item(56, 13).Value = -8362.5  # M56: was -9506.666999999999

$ws.Cells.Item(113, 8).Value = 625.4681  # H113: was 578.69696
$ws.Cells.Item(113, 9).Value = 582.44446  # I113: was 572.3570999999999
$ws.Cells.Item(113, 10).Value = 652.1724  # J113: was 583.3684
$ws.Cells.Item(113, 11).Value = 1747.33338  # K113: was 1717.0713
$ws.Cells.Item(113, 12).Value = 1956.5172  # L113: was 1750.1052
$ws.Cells.Item(113, 13).Value = 422.66662  # M113: was 452.9287000000002
$ws.Cells.Item(113, 14).Value = -6296.5172  # N113: was -6090.1052

$ws.Cells.Item(117, 8).Value = 3550.4211  # H117: was 4467.533
$ws.Cells.Item(117, 10).Value = 3714.3333  # J117: was 4743.7856
$ws.Cells.Item(117, 12).Value = 11142.9999  # L117: was 14231.3568
$ws.Cells.Item(117, 14).Value = -18026.9999  # N117: was -21115.3568

$ws.Cells.Item(121, 8).Value = 1773.9333  # H121: was 1799.5593
$ws.Cells.Item(121, 10).Value = 1773.9333  # J121: was 1799.5593
$ws.Cells.Item(121, 12).Value = 5321.7999  # L121: was 5398.6779
$ws.Cells.Item(121, 14).Value = -7941.7999  # N121: was -8018.6779

$ws.Cells.Item(131, 8).Value = 11364556  # H131: was 7693163.5
$ws.Cells.Item(131, 9).Value = 125000696  # I131: was 166667470
$ws.Cells.Item(131, 10).Value = 941.675  # J131: was 858.2742
$ws.Cells.Item(131, 11).Value = 375002088  # K131: was 500002410
$ws.Cells.Item(131, 12).Value = 2825.025  # L131: was 2574.8226
$ws.Cells.Item(131, 13).Value = -374997048  # M131: was -499997370
$ws.Cells.Item(131, 14).Value = -12905.025  # N131: was -12654.8226

$ws.Cells.Item(132, 8).Value = 1361.6666  # H132: was 1490.6129
$ws.Cells.Item(132, 9).Value = 707.6842  # I132: was 768.55554
$ws.Cells.Item(132, 10).Value = 2092.5881  # J132: was 2490.3845
$ws.Cells.Item(132, 11).Value = 6369.1578  # K132: was 6916.99986
$ws.Cells.Item(132, 12).Value = 18833.2929  # L132: was 22413.4605
$ws.Cells.Item(132, 13).Value = -3839.1578  # M132: was -4386.99986
$ws.Cells.Item(132, 14).Value = -23893.2929  # N132: was -27473.4605

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 34998.5  # H35: was 35000
$ws.Cells.Item(35, 10).Value = 34998.5  # J35: was 35000
$ws.Cells.Item(35, 12).Value = 34998.5  # L35: was 35000
$ws.Cells.Item(35, 14).Value = -35594.5  # N35: was -35596

$ws.Cells.Item(137, 8).Value = 43076  # H137: was 43080
$ws.Cells.Item(137, 10).Value = 43076  # J137: was 43080
$ws.Cells.Item(137, 12).Value = 43076  # L137: was 43080
$ws.Cells.Item(137, 14).Value = -53276  # N137: was -53280

$ws.Cells.Item(140, 8).Value = 38711.25  # H140: was 39933.332
$ws.Cells.Item(140, 10).Value = 38711.25  # J140: was 39933.332
$ws.Cells.Item(140, 12).Value = 38711.25  # L140: was 39933.332
$ws.Cells.Item(140, 14).Value = -49071.25  # N140: was -50293.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2699.6597  # H132: was 2693.791
$ws.Cells.Item(132, 9).Value = 974.17145  # I132: was 911.8095
$ws.Cells.Item(132, 10).Value = 7732.3335  # J132: was 5687.52
$ws.Cells.Item(132, 11).Value = 2922.51435  # K132: was 2735.4285
$ws.Cells.Item(132, 12).Value = 23197.0005  # L132: was 17062.56
$ws.Cells.Item(132, 13).Value = -392.5143500000004  # M132: was -205.4285
$ws.Cells.Item(132, 14).Value = -28257.0005  # N132: was -22122.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 31933.334  # H49: was 30999.75
$ws.Cells.Item(49, 9).Value = 0  # I49: was 10000
$ws.Cells.Item(49, 10).Value = 31933.334  # J49: was 37999.668
$ws.Cells.Item(49, 11).Value = 0  # K49: was 10000
$ws.Cells.Item(49, 12).Value = 31933.334  # L49: was 37999.668
$ws.Cells.Item(49, 13).ClearContents()  # M49: was -9770
$ws.Cells.Item(49, 14).Value = -32393.334  # N49: was -38459.668

$ws.Cells.Item(126, 8).Value = 324666.12  # H126: was 297737.94
$ws.Cells.Item(126, 9).Value = 1667.1765  # I126: was 1415.4286
$ws.Cells.Item(126, 10).Value = 667852.5  # J126: was 712589.4399999999
$ws.Cells.Item(126, 11).Value = 5001.529500000001  # K126: was 4246.2858
$ws.Cells.Item(126, 12).Value = 2003557.5  # L126: was 2137768.32
$ws.Cells.Item(126, 13).Value = -2531.529500000001  # M126: was -1776.2858
$ws.Cells.Item(126, 14).Value = -2008497.5  # N126: was -2142708.32
